# Updates the cryptos list with new price/volume figures, and swaps the
# ranking positions of MXToken and Hedera (rows 38/39).
#
# The "Price" column (D) frequently holds values that look numeric
# (e.g. "1.013", "90.50", "27.437.37") but must be preserved verbatim as
# text (trailing zeros, multi-dot thousands separators, etc.), so the
# cell's NumberFormat is forced to Text ("@") before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $d, $e) {
    if ($null -ne $d) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $d
    }
    $ws.Range("E$row").Value = "  $e  "
}

Set-Row 2  "27.437.37"      "+1.46%"
Set-Row 3  "1.865.25"       "+0.65%"
Set-Row 4  "1.013"          "-0.11%"
Set-Row 5  "311.45"         "+0.32%"
Set-Row 6  "1.011"          "-0.25%"
Set-Row 7  "0.4779"         "-0.04%"
Set-Row 8  "0.3762"         "+2.16%"
Set-Row 9  "0.07327"        "+1.32%"
Set-Row 10 "0.9349"         "+0.34%"
Set-Row 11 "20.67"          "+4.66%"
Set-Row 12 "0.07827"        "+1.18%"
Set-Row 13 "1.883.69"       "+1.25%"
Set-Row 14 "5.435"          "+1.94%"
Set-Row 15 "6.552"          "+1.78%"
Set-Row 16 "90.50"          "+1.65%"
Set-Row 17 "1.014"          "-0.20%"
Set-Row 18 "0.000008892"    "+2.84%"
Set-Row 19 $null             "-0.15%"
Set-Row 20 "27.496.17"      "+1.59%"
Set-Row 21 "14.73"          "+1.20%"
Set-Row 22 "5.115"          "+1.03%"
Set-Row 23 "10.69"          "+0.09%"
Set-Row 24 "1.938"          "+0.24%"
Set-Row 25 "155.34"         "+1.58%"
Set-Row 26 "18.47"          "+1.31%"
Set-Row 27 $null             "+0.75%"
Set-Row 28 "115.46"         "+0.82%"
Set-Row 29 "4.964"          "-0.69%"
Set-Row 30 "0.08898"        "-0.08%"
Set-Row 31 "3.328"          "-0.09%"
Set-Row 32 "1.214"          "+3.01%"
Set-Row 33 "0.7565"         "+1.42%"
Set-Row 34 "4.600"          "+2.03%"
Set-Row 35 "2.736"          "+0.10%"
Set-Row 36 $null             "+0.77%"
Set-Row 37 "0.02034"        "+3.98%"

# Row 38/39: MXToken and Hedera swap ranking positions.
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-Row 38 "2.990"          "+0.39%"

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-Row 39 "0.05257"        "-0.44%"

Set-Row 40 "0.5312"         "+1.81%"
Set-Row 41 "7.079"          "+0.67%"
Set-Row 42 "8.567"          "+4.09%"
Set-Row 43 "0.1524"         "+0.79%"
Set-Row 44 "10.64"          "+0.56%"
Set-Row 45 "0.4804"         "+1.20%"
Set-Row 46 "1.011"          "-0.47%"
Set-Row 47 "1.656"          "+2.78%"
Set-Row 48 "102.80"         "+0.96%"
Set-Row 49 "67.32"          "+1.99%"
Set-Row 51 "0.9176"         "+3.16%"
